$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (238) down
# across the six new rows (239-244) so column A keeps its date style (s="2").
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 239: 2021-04-27
$ws.Range("A239").Value = 44313
$ws.Range("B239").Value = 0
$ws.Range("C239").Value = 10
$ws.Range("D239").Value = 238.2654276864427

# Row 240: 2021-04-28
$ws.Range("A240").Value = 44314
$ws.Range("B240").Value = 0
$ws.Range("C240").Value = 10
$ws.Range("D240").Value = 238.2654276864427

# Row 241: 2021-04-29
$ws.Range("A241").Value = 44315
$ws.Range("B241").Value = 1
$ws.Range("C241").Value = 10
$ws.Range("D241").Value = 238.2654276864427

# Row 242: 2021-04-30
$ws.Range("A242").Value = 44316
$ws.Range("B242").Value = 2
$ws.Range("C242").Value = 10
$ws.Range("D242").Value = 238.2654276864427

# Row 243: 2021-05-01
$ws.Range("A243").Value = 44317
$ws.Range("B243").Value = 0
$ws.Range("C243").Value = 9
$ws.Range("D243").Value = 214.4388849177984

# Row 244: 2021-05-02
$ws.Range("A244").Value = 44318
$ws.Range("B244").Value = 0
$ws.Range("C244").Value = 3
$ws.Range("D244").Value = 71.47962830593281
